# Add 20 more "presence" columns (EC:EV) to sheet1, continuing the existing
# repeating Alain/Henri/Tony/Dulcinee header pattern (and the corresponding
# OUI/NON answer pattern in the data rows). This pushes the trailing
# email/status columns that used to live at EC:ED out to EW:EX.
#
# The new columns repeat the same 4-column cycle already used across the
# sheet. Columns I:AB are one such 20-column-wide (5 x 4-cycle) block, so we
# copy it and use "Insert Copied Cells" (Copy, then Insert with
# xlShiftToRight) at EC so the existing EC:ED data is shifted right to
# EW:EX and the new cells land with the same formatting as their source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Range("I1:AB9").Copy()
$ws.Range("EC1:EV9").Insert(-4161)

$wb.Application.CutCopyMode = 0
